$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Extend the table: copy the bordered/wrap-text formatting of the
#    existing data rows (2-8) down onto the new rows (9-18) first, so
#    the new rows pick up the same cell style (border + wrap) used by
#    the rest of the table.
# ------------------------------------------------------------------
$ws.Range("A2:G8").Copy() | Out-Null
$ws.Range("A9:G18").PasteSpecial(-4122) | Out-Null

# ------------------------------------------------------------------
# 2. (Re)write the test-case table data, rows 2-18, columns A-G.
#    Rows 2-8 keep their S/N, Description, Test Inputs, Expected
#    Results and Actual Results, but the "Test Procedure" wording is
#    updated to mention selecting a role first. Rows 9-18 are brand
#    new admin-login / role-mismatch test cases.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1.1
$ws.Range("C2").Value = "Validate that student with correct credentials can login"
$ws.Range("D2").Value = "Username: amy.ng.2009`nPassword: qwerty128"
$ws.Range("E2").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F2").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G2").Value = "Login successful but redirected to 'plan bid' instead of 'home' page"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1.1
$ws.Range("C3").Value = "Validate that student who did not key in username would not be able to login successfully and show error message"
$ws.Range("D3").Value = "Username: `nPassword: qwerty128"
$ws.Range("E3").Value = "Select student as their role, do not key in anything in the username input field, put password into password input field. Submit"
$ws.Range("F3").Value = "Login failed, show error message"
$ws.Range("G3").Value = "Matched expected results"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1.1
$ws.Range("C4").Value = "Validate that student who did not key in password would not be able to login successfully and show error message"
$ws.Range("D4").Value = "Username: amy.ng.2009`nPassword: "
$ws.Range("E4").Value = "Select student as their role, put username into username input field, do not key in anything in the password input field. Submit"
$ws.Range("F4").Value = "Login failed, show error message"
$ws.Range("G4").Value = "Matched expected results"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1.1
$ws.Range("C5").Value = "Validate that student with incorrect username would not be able to login and show error message"
$ws.Range("D5").Value = "Username: amy.ng`nPassword: qwerty128"
$ws.Range("E5").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F5").Value = "Login failed, show error message"
$ws.Range("G5").Value = "Matched expected results"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1.1
$ws.Range("C6").Value = "Validate that student with incorrect password would not be able to login and show error message"
$ws.Range("D6").Value = "Username: amy.ng.2009`nPassword: 123456"
$ws.Range("E6").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F6").Value = "Login failed, show error message"
$ws.Range("G6").Value = "Matched expected results"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 1.1
$ws.Range("C7").Value = "Validate that student with correct credentials can login (double check)"
$ws.Range("D7").Value = "Username: ben.ng.2009`nPassword: qwerty129"
$ws.Range("E7").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F7").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G7").Value = "Login successful but wrong name was displayed on home page"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1.1
$ws.Range("C8").Value = "Validate that student with correct credentials can login (double double check)"
$ws.Range("D8").Value = "Username: calvin.ng.2009`nPassword: qwerty130"
$ws.Range("E8").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F8").Value = "Login success, student home page displayed with correct bidding summary"
$ws.Range("G8").Value = "Login successful, bidding summary matches expected result"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 1.1
$ws.Range("C9").Value = "Validate that admin with correct credentials can login"
$ws.Range("D9").Value = "Username: admin`nPassword: adminpassword"
$ws.Range("E9").Value = "Select admin as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F9").Value = "Login success, bidding details for admin is displayed."
$ws.Range("G9").Value = "Login unsuccessful and error displayed (got directed to student home page and hence causing error as admin do not have correct bidding summary)"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 1.1
$ws.Range("C10").Value = "Validate that admin with correct credentials can login (double check)"
$ws.Range("D10").Value = "Username: admin`nPassword: adminpassword"
$ws.Range("E10").Value = "Select admin as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F10").Value = "Login success, bidding details for admin is displayed."
$ws.Range("G10").Value = "Login successful, bidding admin page matches expected result"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 1.1
$ws.Range("C11").Value = "Validate that if admin did not key in username, they would not be able to login successfully and show error message"
$ws.Range("D11").Value = "Username: `nPassword: adminpassword"
$ws.Range("E11").Value = "Select admin as their role, do not key in anything in the username input field, put password into password input field. Submit"
$ws.Range("F11").Value = "Login failed, show error message"
$ws.Range("G11").Value = "Matched expected results"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 1.1
$ws.Range("C12").Value = "Validate that if admin did not key in password, they would not be able to login successfully and show error message"
$ws.Range("D12").Value = "Username: admin`nPassword: "
$ws.Range("E12").Value = "Select admin as their role, put username into username input field, do not key in anything in the password input field. Submit"
$ws.Range("F12").Value = "Login failed, show error message"
$ws.Range("G12").Value = "Matched expected results"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 1.1
$ws.Range("C13").Value = "Validate that admin with incorrect username would not be able to login and show error message"
$ws.Range("D13").Value = "Username: admin123`nPassword: password"
$ws.Range("E13").Value = "Select admin as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F13").Value = "Login failed, show error message"
$ws.Range("G13").Value = "Matched expected results"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 1.1
$ws.Range("C14").Value = "Validate that admin with incorrect password would not be able to login and show error message"
$ws.Range("D14").Value = "Username: admin`nPassword: woshiadmin"
$ws.Range("E14").Value = "Select admin as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F14").Value = "Login failed, show error message"
$ws.Range("G14").Value = "Matched expected results"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 1.1
$ws.Range("C15").Value = "Ensure that admin would not be able to login if they key in their password in CAPS "
$ws.Range("D15").Value = "Username: admin`nPassword: ADMINPASSWORD"
$ws.Range("E15").Value = "Select admin as their role, put username into username input field, key in password in CAPS into password input field. Submit"
$ws.Range("F15").Value = "Login failed, show error message"
$ws.Range("G15").Value = "Login successful and was brought to the bidding admin page"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 1.1
$ws.Range("C16").Value = "Ensure that admin would not be able to login if they key in their password in CAPS - double check"
$ws.Range("D16").Value = "Username: admin`nPassword: ADMINPASSWORD"
$ws.Range("E16").Value = "Select admin as their role, put username into username input field, key in password in CAPS into password input field. Submit"
$ws.Range("F16").Value = "Login failed, show error message"
$ws.Range("G16").Value = "Matched expected results"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 1.1
$ws.Range("C17").Value = "Validate if student is able to login into admin home page using student's credentials"
$ws.Range("D17").Value = "Username: amy.ng.2009`nPassword: qwerty128"
$ws.Range("E17").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F17").Value = "Login failed, show error message"
$ws.Range("G17").Value = "Login successful and was brought to the bidding admin page"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 1.1
$ws.Range("C18").Value = "Validate if student is able to login into admin home page using student's credentials"
$ws.Range("D18").Value = "Username: amy.ng.2009`nPassword: qwerty128"
$ws.Range("E18").Value = "Select student as their role, put username into username input field, put password into password input field. Submit"
$ws.Range("F18").Value = "Login failed, show error message"
$ws.Range("G18").Value = "Matched expected results"

# ------------------------------------------------------------------
# 3. Highlight the "Actual Results" column: red for mismatched /
#    failing cases, green for cases that matched the expected
#    result. Color the first red cell and first green cell before
#    anything else so the red fill is registered first.
# ------------------------------------------------------------------
$ws.Range("G2").Interior.Color = 8487423
$ws.Range("G3").Interior.Color = 6280092
$ws.Range("G4").Interior.Color = 6280092
$ws.Range("G5").Interior.Color = 6280092
$ws.Range("G6").Interior.Color = 6280092
$ws.Range("G7").Interior.Color = 8487423
$ws.Range("G8").Interior.Color = 6280092
$ws.Range("G9").Interior.Color = 8487423
$ws.Range("G10").Interior.Color = 6280092
$ws.Range("G11").Interior.Color = 6280092
$ws.Range("G12").Interior.Color = 6280092
$ws.Range("G13").Interior.Color = 6280092
$ws.Range("G14").Interior.Color = 6280092
$ws.Range("G15").Interior.Color = 8487423
$ws.Range("G16").Interior.Color = 6280092
$ws.Range("G17").Interior.Color = 8487423
$ws.Range("G18").Interior.Color = 6280092

# ------------------------------------------------------------------
# 4. Row heights: the wrapped text in the new/edited rows needs more
#    vertical space than the sheet default.
# ------------------------------------------------------------------
$ws.Rows(2).RowHeight = 43.2
$ws.Rows(3).RowHeight = 43.2
$ws.Rows(4).RowHeight = 43.2
$ws.Rows(5).RowHeight = 43.2
$ws.Rows(6).RowHeight = 43.2
$ws.Rows(7).RowHeight = 43.2
$ws.Rows(8).RowHeight = 43.2
$ws.Rows(9).RowHeight = 57.6
$ws.Rows(10).RowHeight = 43.2
$ws.Rows(11).RowHeight = 43.2
$ws.Rows(12).RowHeight = 43.2
$ws.Rows(13).RowHeight = 43.2
$ws.Rows(14).RowHeight = 43.2
$ws.Rows(15).RowHeight = 43.2
$ws.Rows(16).RowHeight = 43.2
$ws.Rows(17).RowHeight = 43.2
$ws.Rows(18).RowHeight = 43.2

# ------------------------------------------------------------------
# 5. Leave the view scrolled/selected near the newly added rows.
# ------------------------------------------------------------------
$ws.Range("E14").Select() | Out-Null
